# Journal_Travail_Guillaume.xlsx — add a new journal entry (row 37):
# "Problème avec les id généré automatiquement dans derby", 3h, same date as
# the two entries above (14/04/2018), plus fill in the previously-blank
# hours (1h) for the "ClientRepository pour derby" entry on row 36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 36 (développement et test de ClientRepository pour derby): was blank,
# now logged as 1 hour.
$ws.Range("C36").Value = 1

# Row 37 was a fully blank placeholder row — fill it in with a new entry
# dated the same day as row 35/36 (14/04/2018, serial 43204), describing the
# new problem encountered and the hours spent on it.
$ws.Range("A37").Value = $ws.Range("A36").Value()
$ws.Range("B37").Value = "Problème avec les id généré automatiquement dans derby"
$ws.Range("C37").Value = 3

# Reflect where the author was last looking/working in the sheet when they
# saved (view state only — no data impact).
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("I36").Select()

$wb.Save()
